$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07908133333333334
$ws.Range("H2").Value = 0.237244
$ws.Range("I2").Value = 0.1380838080781507
$ws.Range("J2").Value = 0.1380838080781507
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.063701
$ws.Range("N2").Value = 0.191103
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 0.005037560014666666
$ws.Range("R2").Value = 0.045338040132
$ws.Range("S2").Value = 0.002224769949777768
$ws.Range("T2").Value = 0.002224769949777767

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07908133333333334
$ws.Range("H3").Value = 0.237244
$ws.Range("I3").Value = 0.1380838080781507
$ws.Range("J3").Value = 0.1380838080781507
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 0.2136379320351111
$ws.Range("R3").Value = 1.922741388316
$ws.Range("S3").Value = 0.0943502906050898
$ws.Range("T3").Value = 0.09435029060508979

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07908133333333334
$ws.Range("H4").Value = 0.237244
$ws.Range("I4").Value = 0.1380838080781507
$ws.Range("J4").Value = 0.1380838080781507
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 0.09398850735244445
$ws.Range("R4").Value = 0.845896566172
$ws.Range("S4").Value = 0.04150874752328316
$ws.Range("T4").Value = 0.04150874752328314

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.07720866666666666
$ws.Range("H5").Value = 0.231626
$ws.Range("I5").Value = 0.1348139473702591
$ws.Range("J5").Value = 0.134813947370259
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.063701
$ws.Range("N5").Value = 0.191103
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("Q5").Value = 0.004918269275333333
$ws.Range("R5").Value = 0.044264423478
$ws.Range("S5").Value = 0.002172086815208079
$ws.Range("T5").Value = 0.002172086815208077

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ccl21b"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.07720866666666666
$ws.Range("H6").Value = 0.231626
$ws.Range("I6").Value = 0.1348139473702591
$ws.Range("J6").Value = 0.134813947370259
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 0.2085789299015555
$ws.Range("R6").Value = 1.877210369114
$ws.Range("S6").Value = 0.09211605103477656
$ws.Range("T6").Value = 0.09211605103477655

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ccl21b"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.07720866666666666
$ws.Range("H7").Value = 0.231626
$ws.Range("I7").Value = 0.1348139473702591
$ws.Range("J7").Value = 0.134813947370259
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 0.09176283490422223
$ws.Range("R7").Value = 0.825865514138
$ws.Range("S7").Value = 0.04052580952027442
$ws.Range("T7").Value = 0.0405258095202744

# Row 8
$ws.Range("A8").Value = "Inflammatory-Mac"
$ws.Range("B8").Value = "Ccl21b"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.03917266666666667
$ws.Range("H8").Value = 0.117518
$ws.Range("I8").Value = 0.06839933974190335
$ws.Range("J8").Value = 0.06839933974190333
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.063701
$ws.Range("N8").Value = 0.191103
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 0.002495338039333333
$ws.Range("R8").Value = 0.022458042354
$ws.Range("S8").Value = 0.001102032148159632
$ws.Range("T8").Value = 0.001102032148159632

# Row 9
$ws.Range("A9").Value = "Inflammatory-Mac"
$ws.Range("B9").Value = "Ccl21b"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.03917266666666667
$ws.Range("H9").Value = 0.117518
$ws.Range("I9").Value = 0.06839933974190335
$ws.Range("J9").Value = 0.06839933974190333
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 0.1058248153668889
$ws.Range("R9").Value = 0.9524233383020001
$ws.Range("S9").Value = 0.04673609217231603
$ws.Range("T9").Value = 0.04673609217231602

# Row 10
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Ccl21b"
$ws.Range("C10").Value = "Cxcr3"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.03917266666666667
$ws.Range("H10").Value = 0.117518
$ws.Range("I10").Value = 0.06839933974190335
$ws.Range("J10").Value = 0.06839933974190333
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 0.04655688408155556
$ws.Range("R10").Value = 0.419011956734
$ws.Range("S10").Value = 0.02056121542142769
$ws.Range("T10").Value = 0.02056121542142768

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Ccl21b"
$ws.Range("C11").Value = "Cxcr3"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3772426666666667
$ws.Range("H11").Value = 1.131728
$ws.Range("I11").Value = 0.658702904809687
$ws.Range("J11").Value = 0.6587029048096869
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.063701
$ws.Range("N11").Value = 0.191103
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.02403073510933333
$ws.Range("R11").Value = 0.216276615984
$ws.Range("S11").Value = 0.01061284772522
$ws.Range("T11").Value = 0.01061284772522

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ccl21b"
$ws.Range("C12").Value = "Cxcr3"
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3772426666666667
$ws.Range("H12").Value = 1.131728
$ws.Range("I12").Value = 0.658702904809687
$ws.Range("J12").Value = 0.6587029048096869
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 1.019119680776889
$ws.Range("R12").Value = 9.172077126992
$ws.Range("S12").Value = 0.4500803631953477
$ws.Range("T12").Value = 0.4500803631953476

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ccl21b"
$ws.Range("C13").Value = "Cxcr3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3772426666666667
$ws.Range("H13").Value = 1.131728
$ws.Range("I13").Value = 0.658702904809687
$ws.Range("J13").Value = 0.6587029048096869
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 0.4483545440515556
$ws.Range("R13").Value = 4.035190896464
$ws.Range("S13").Value = 0.1980096938891192
$ws.Range("T13").Value = 0.1980096938891192

